# Update the handback status report timestamps (simulating a fresh report generation).
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first data row.
$wsOverview.Range("G2").Value = "2016-08-21 19:07:04"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the first data row.
$wsZhCn.Range("H2").Value = "2016-08-21 19:06:57"
$wsZhCn.Range("K2").Value = "2016-08-21 19:07:26"

# de-de sheet: "Correspond Handback DateTime" (K) for the first data row.
$wsDeDe.Range("K2").Value = "2016-08-21 19:07:33"
